# Budget tracker update: fix encoding on an existing entry, replace the
# "onlyfans" row with a new expense, log a new "HYRE" expense, and start
# tracking April in its own sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("March")

# --- Fix mangled characters in the "AtB månedskort" transport entry (row 6) ---
$ws.Range("B6").Value = "AtB mÃ¥nedskort"

# --- Row 7: replace the old "onlyfans" entry with a new Food expense ---
$ws.Range("A7").Value = "Food"
$ws.Range("B7").Value = "asdfg"

$ws.Range("C7").Value = "2023-03-15"
$ws.Range("C7").ClearFormats()

$ws.Range("D7").Value = "1000.0"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "Checkings"

# --- Row 8: new Transportation expense (HYRE) ---
$ws.Range("A8").Value = "Transportation"
$ws.Range("B8").Value = "HYRE"

$ws.Range("C8").Value = "2023-03-16"
$ws.Range("C8").ClearFormats()

$ws.Range("D8").Value = "150.0"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "Card"

# --- Add the "April" sheet right after "March" ---
$april = $wb.Worksheets.Add()
$april.Name = "April"
$april.Move($null, $ws)

$april.Range("A1").Value = "Category"
$april.Range("B1").Value = "Name"
$april.Range("C1").Value = "Date"
$april.Range("D1").Value = "Price"
$april.Range("E1").Value = "Account"

$april.Range("A2").Value = "Clothing"
$april.Range("B2").Value = "gucci sokker"

$april.Range("C2").Value = "2023-04-06"
$april.Range("C2").ClearFormats()

$april.Range("D2").Value = "2000.0"
$april.Range("D2").ClearFormats()

$april.Range("E2").Value = "Card"
